# report-template-hw3.docx cleanup:
#  - remove the stray _GoBack bookmark near the top and re-add it at the
#    very end of the body (after the last {@hw_9} run)
#  - merge runs that were split apart by Word's spell/grammar checker
#    (the proofErr spellStart/spellEnd/gramStart/gramEnd markers go away
#    automatically once the surrounding text is re-written as a single run)
#  - drop the four trailing empty paragraphs at the end of the body

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# ---- body: merge runs split by proofErr tags --------------------------
Replace-Text "{name}" "{name}"
Replace-Text "Matr.-Nr.:" "Matr.-Nr.:"
Replace-Text "{matrikelnummer}" "{matrikelnummer}"
Replace-Text "Megenoperationen" "Megenoperationen"
Replace-Text "{group} " "{group} "
Replace-Text "{date} {start} – {end}" "{date} {start} – {end}"
Replace-Text "{passes}/{all}" "{passes}/{all}"
Replace-Text "{#hasFailures}" "{#hasFailures}"
Replace-Text "{#failureFiles}" "{#failureFiles}"
Replace-Text "{filename}:" "{filename}:"
Replace-Text " {message}" " {message}"
Replace-Text "{/failureFiles}{/hasFailures}" "{/failureFiles}{/hasFailures}"
Replace-Text "Fertigungsdaten VIEW: {hw_2_file}" "Fertigungsdaten VIEW: {hw_2_file}"

# ---- header: merge runs split by proofErr tags -------------------------
Replace-Text "{name} " "{name} "
Replace-Text "3. Hausübung" "3. Hausübung"
Replace-Text "{matrikelnummer}  " "{matrikelnummer}  "

Write-Output "done"
